$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from A11 (bold/centered/bordered label style) to A12
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 114
$ws.Range("C12").Value = 808
$ws.Range("D12").Value = 258
$ws.Range("E12").Value = "'"
$ws.Range("F12").Value = "'"
$ws.Range("G12").Value = 68
$ws.Range("H12").Value = 2217
$ws.Range("I12").Value = 348
$ws.Range("J12").Value = 2418
$ws.Range("K12").Value = "'"
$ws.Range("L12").Value = 98
$ws.Range("M12").Value = 593
$ws.Range("N12").Value = 82
$ws.Range("O12").Value = "'"
$ws.Range("P12").Value = 132
$ws.Range("Q12").Value = 79
$ws.Range("R12").Value = 53
$ws.Range("S12").Value = 120
$ws.Range("T12").Value = 53
$ws.Range("U12").Value = 151
$ws.Range("V12").Value = 351
$ws.Range("W12").Value = 3
$ws.Range("X12").Value = 761
$ws.Range("Y12").Value = 3611
$ws.Range("Z12").Value = 3152
$ws.Range("AA12").Value = 752
$ws.Range("AB12").Value = "'"
$ws.Range("AC12").Value = 100
$ws.Range("AD12").Value = 16
$ws.Range("AE12").Value = 44
$ws.Range("AF12").Value = 1593
$ws.Range("AG12").Value = 1951
$ws.Range("AH12").Value = 465
$ws.Range("AI12").Value = 15
$ws.Range("AJ12").Value = 138
$ws.Range("AK12").Value = 36
$ws.Range("AL12").Value = "'"
$ws.Range("AM12").Value = "'"
$ws.Range("AN12").Value = 257
$ws.Range("AO12").Value = 169
$ws.Range("AP12").Value = 301
$ws.Range("AQ12").Value = "'"
$ws.Range("AR12").Value = 10
$ws.Range("AS12").Value = 3166
$ws.Range("AT12").Value = 44
$ws.Range("AU12").Value = 7
$ws.Range("AV12").Value = 133
$ws.Range("AW12").Value = 1238
$ws.Range("AX12").Value = 521
$ws.Range("AY12").Value = 269
$ws.Range("AZ12").Value = 118
$ws.Range("BA12").Value = "'"
$ws.Range("BB12").Value = 287
$ws.Range("BC12").Value = 580
$ws.Range("BD12").Value = 396
$ws.Range("BE12").Value = 2922
$ws.Range("BF12").Value = "'"
$ws.Range("BG12").Value = 32036
$ws.Range("BH12").Value = 105
$ws.Range("BI12").Value = "'"
$ws.Range("BJ12").Value = 745
$ws.Range("BK12").Value = 188
